$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Exact "English" "ภาษาอังกฤษ"
Replace-Exact " / Portuguese / French / Thai / Vietnamese / Spanish" " / ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน"
Replace-Exact "English" "ภาษาอังกฤษ"

Replace-Exact "Brief" "บทย่อ"
Replace-Exact "An email sent to partners in the target country who have sent their documents for review. It will be sent via customer.io" "อีเมลที่ส่งถึงพันธมิตรในประเทศเป้าหมายที่ได้ส่งเอกสารของพวกเขาสำหรับการตรวจสอบ โดยมันจะถูกส่งผ่านทาง customer.io"
Replace-Exact "Target audience" "กลุ่มเป้าหมาย"
Replace-Exact "Invited partners who have submitted their documents" "พันธมิตรที่ได้รับเชิญซึ่งได้ส่งเอกสารของพวกเขาแล้ว"

Replace-Exact "Subject line" "หัวเรื่อง"
Replace-Exact " — we got your docs!  " " — เราได้รับเอกสารของคุณแล้ว!  "

Replace-Exact "Thank you for submitting your documents" "ขอบคุณที่ส่งเอกสารของคุณมาให้เรา"

Replace-Exact "Hi " "สวัสดี "
Replace-Exact ", " " "

Replace-Exact "Thank you for providing us with your documents for the upcoming " "ขอขอบคุณที่ส่งเอกสารของคุณให้กับเราสำหรับงาน "
Replace-Exact ". Based on the information you’ve given us, we’ll make the necessary arrangements, including accommodation and transportation." " ที่กำลังจะเกิดขึ้น จากข้อมูลที่คุณได้ให้ไว้กับเรา เราจะดำเนินการด้านต่างๆ ตามที่จำเป็นรวมถึงจัดการเรื่องที่พักและการเดินทาง"

Replace-Exact "We’re currently reviewing your documents and will reach out to you if we need anything else. " "ขณะนี้ เรากำลังตรวจสอบเอกสารของคุณ และจะติดต่อหาคุณหากเราต้องการสิ่งอื่นใดเพิ่มเติม "

Replace-Exact "If you have any questions, please contact us via " "หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง "
Replace-Exact "live chat" "แชทสด"
Replace-Exact " or " " หรือทาง "
Replace-Exact ". " " "

Replace-Exact "If you have any questions, please contact your country manager, " "หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ "
Replace-Exact ", at " " ที่ "
Replace-Exact " or " " หรือ "
Replace-Exact " (WhatsApp). " " (WhatsApp) "

Replace-Exact "We look forward to seeing you at " "เราหวังเป็นอย่างยิ่งว่าจะได้พบคุณที่ "
Replace-Exact ". " " "

Replace-Exact "choose either one" "เลือกอย่างใดอย่างหนึ่ง"
